$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "App width: $($excel.Width) height: $($excel.Height)"
$excel.Width = 30360
$excel.Height = 19540
Write-Host "App width: $($excel.Width) height: $($excel.Height)"
